$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The C10 cell (row "R30") currently holds 18 and should be updated to 100,
# per the authoritative diff of the workbook's commit.
$ws.Range("C10").Value = 100
